# Fruta / hortaliza, semanal
# Insert a new weekly record row above the current row 140 (Granada / Vega Modelo de Temuco),
# pushing the existing rows 140-192 down to 141-193, and fill the new row with the
# latest week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 140; existing rows 140:192 shift down to 141:193
$ws.Rows.Item(140).Insert()

# Populate the new row 140 with the new weekly observation
$ws.Cells.Item(140, 1).Value = 10                                 # A Mercado ID
$ws.Cells.Item(140, 2).Value = "Vega Modelo de Temuco"             # B Mercado
$ws.Cells.Item(140, 3).Value = "La Araucanía"                      # C Región
$ws.Cells.Item(140, 4).Value = 44875                                # D Fecha
$ws.Cells.Item(140, 5).Value = 9                                    # E Codreg
$ws.Cells.Item(140, 6).Value = "Fruta"                              # F Tipo
$ws.Cells.Item(140, 7).Value = 100104                               # G Producto ID
$ws.Cells.Item(140, 8).Value = "Frutos de pepita"                   # H Producto
$ws.Cells.Item(140, 9).Value = 100104001                            # I Categoría ID
$ws.Cells.Item(140, 10).Value = "Granada"                           # J Categoría
$ws.Cells.Item(140, 11).Value = "Wonderfull"                        # K Variedad
$ws.Cells.Item(140, 12).Value = "Primera"                           # L Calidad
$ws.Cells.Item(140, 13).Value = 100                                 # M Volumen
$ws.Cells.Item(140, 14).Value = 16000                               # N Precio mínimo
$ws.Cells.Item(140, 15).Value = 16000                               # O Precio máximo
$ws.Cells.Item(140, 16).Value = 16000                               # P Precio promedio ponderado
$ws.Cells.Item(140, 17).Value = "$/bandeja 15 kilos granel"         # Q Unidad de comercialización
$ws.Cells.Item(140, 18).Value = "Provincia de Limarí"                # R Origen
$ws.Cells.Item(140, 19).Value = 1067                                 # S Precio $/Kg
$ws.Cells.Item(140, 20).Value = 15                                   # T Kg / unidad
